$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "247.58"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.91"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.394"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05639"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.432"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.369"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8183"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9358"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01155"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1440"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07504"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03256"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09320"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.559"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001596"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04736"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006364"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005060"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001036"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001501"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.754"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.188"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3307"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1328"
$ws.Range("E26").Value = "25ProBitTokenPROB"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0003002"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03951"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006934"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1065"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003402"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008519"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005584"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005502"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.7805"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1776"
$ws.Range("E49").Value = "48BOLOBOLO"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.01011"
